$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New column H ("2022") mirrors the existing G ("2021") column: same section
# layout, same borders, same fonts -- just one more year of data with a
# one-decimal ("0.0") number format.
# ---------------------------------------------------------------------------

# Row 4 is the year-header row; H4 should look exactly like G4.
$ws.Range("G4").Copy() | Out-Null
$ws.Range("H4").PasteSpecial(-4122) | Out-Null
$ws.Range("H4").Value = 2022

# Row 5 ("Кыргыз Республикасы" / national total) got its own one-off style
# in the original edit: bold 9pt "Times New Roman Cyr" (same font already
# used for the bold Cyrillic header cells), 0.0 number format, no border,
# no explicit alignment override.
$ws.Range("B4").Copy() | Out-Null
$ws.Range("H5").PasteSpecial(-4122) | Out-Null
$ws.Range("H5").Value = 5.6504815716081698
$ws.Range("H5").NumberFormat = "0.0"
$ws.Range("H5").Borders.LineStyle = -4142
$ws.Range("H5").HorizontalAlignment = 1
$ws.Range("H5").VerticalAlignment = -4107

# Rows 6-36: same regular (non-bold) 9pt "Times New Roman" font as column G,
# right/center aligned, no border -- just add the 0.0 number format.
$values = @{
    7  = 11.990552717185041
    8  = 0.36686386493060885
    10 = 4.7585078053197183
    11 = 6.2147211180387529
    13 = 6.6299725226642234
    14 = 1.0045350275012754
    15 = 15.32109744080277
    16 = 11.932654937916501
    17 = 4.4024960942722968
    18 = 13.275841712798133
    19 = 5.5953228746387378
    20 = 5.7561942305949083
    21 = 0.91260128840317045
    23 = 0.12558892880771302
    24 = 6.813728136595028
    25 = 7.4835121062312364
    27 = 1.0698262411858328
    28 = 2.525977374670846
    29 = 6.9014261042903025
    30 = 7.9091356334900151
    31 = 3.3800067710254136
    33 = 4.7357406632935053
    34 = 4.7664658340238164
    35 = 5.4209310439574798
    36 = 6.4917222807546029
}

foreach ($r in 6..36) {
    $ws.Range("G$r").Copy() | Out-Null
    $ws.Range("H$r").PasteSpecial(-4122) | Out-Null
    $ws.Range("H$r").NumberFormat = "0.0"
}

foreach ($r in $values.Keys) {
    $ws.Range("H$r").Value = $values[$r]
}

# Row 37 is the last data row (has the bottom medium border under the table).
$ws.Range("G37").Copy() | Out-Null
$ws.Range("H37").PasteSpecial(-4122) | Out-Null
$ws.Range("H37").Value = 6.4231110817165149
$ws.Range("H37").NumberFormat = "0.0"

# ---------------------------------------------------------------------------
# View state: scroll back to the top of the sheet and leave the selection on
# D1 (matches the saved sheetView in the workbook after the edit).
# ---------------------------------------------------------------------------
$ws.Range("D1").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
